# Update "想去人数" (interested-attendee count) values in column F.
# Sheet 1 = "展览" (Exhibitions) and Sheet 4 = "全部类型" (All Types) both list
# (a subset of) the same events, so the refreshed counts are applied to both.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 ---
$updates1 = @(
    @{ Row = 2;  Old = 60;    New = 64 }
    @{ Row = 3;  Old = 21599; New = 21623 }
    @{ Row = 4;  Old = 820;   New = 821 }
    @{ Row = 5;  Old = 341;   New = 340 }
    @{ Row = 8;  Old = 8032;  New = 8041 }
    @{ Row = 9;  Old = 561;   New = 564 }
    @{ Row = 10; Old = 47;    New = 49 }
    @{ Row = 11; Old = 768;   New = 770 }
    @{ Row = 12; Old = 323;   New = 326 }
    @{ Row = 15; Old = 185;   New = 187 }
    @{ Row = 20; Old = 562;   New = 564 }
    @{ Row = 21; Old = 86;    New = 87 }
    @{ Row = 23; Old = 56;    New = 57 }
    @{ Row = 24; Old = 93;    New = 94 }
    @{ Row = 25; Old = 89;    New = 90 }
    @{ Row = 26; Old = 359;   New = 362 }
    @{ Row = 27; Old = 1206;  New = 1210 }
    @{ Row = 28; Old = 65;    New = 66 }
    @{ Row = 29; Old = 44;    New = 45 }
    @{ Row = 30; Old = 233;   New = 234 }
    @{ Row = 31; Old = 609;   New = 610 }
    @{ Row = 32; Old = 11;    New = 12 }
    @{ Row = 33; Old = 155;   New = 156 }
    @{ Row = 34; Old = 5143;  New = 5155 }
    @{ Row = 35; Old = 37;    New = 38 }
    @{ Row = 37; Old = 59;    New = 61 }
    @{ Row = 39; Old = 13277; New = 13298 }
    @{ Row = 40; Old = 1374;  New = 1375 }
    @{ Row = 41; Old = 149;   New = 151 }
    @{ Row = 44; Old = 330;   New = 331 }
    @{ Row = 45; Old = 456;   New = 459 }
    @{ Row = 46; Old = 4075;  New = 4077 }
    @{ Row = 47; Old = 29;    New = 31 }
)

foreach ($u in $updates1) {
    $cell = $ws1.Range("F" + $u.Row)
    $actual = $cell.Value2
    if ($actual -ne $u.Old) {
        Write-Output ("WARNING: 展览 F" + $u.Row + " expected " + $u.Old + " but found " + $actual)
    }
    $cell.Value = $u.New
}

# --- Sheet 4: 全部类型 ---
$updates4 = @(
    @{ Row = 2;  Old = 60;    New = 64 }
    @{ Row = 3;  Old = 21599; New = 21623 }
    @{ Row = 6;  Old = 8032;  New = 8041 }
    @{ Row = 7;  Old = 561;   New = 564 }
    @{ Row = 8;  Old = 47;    New = 49 }
    @{ Row = 9;  Old = 768;   New = 770 }
    @{ Row = 10; Old = 323;   New = 326 }
    @{ Row = 13; Old = 185;   New = 187 }
    @{ Row = 18; Old = 562;   New = 564 }
    @{ Row = 19; Old = 86;    New = 87 }
    @{ Row = 21; Old = 56;    New = 57 }
    @{ Row = 22; Old = 93;    New = 94 }
    @{ Row = 23; Old = 89;    New = 90 }
    @{ Row = 24; Old = 359;   New = 362 }
    @{ Row = 25; Old = 1206;  New = 1210 }
    @{ Row = 26; Old = 65;    New = 66 }
    @{ Row = 27; Old = 44;    New = 45 }
    @{ Row = 28; Old = 233;   New = 234 }
    @{ Row = 30; Old = 609;   New = 610 }
    @{ Row = 32; Old = 11;    New = 12 }
    @{ Row = 33; Old = 155;   New = 156 }
    @{ Row = 35; Old = 5143;  New = 5155 }
    @{ Row = 36; Old = 37;    New = 38 }
    @{ Row = 38; Old = 59;    New = 61 }
    @{ Row = 40; Old = 13277; New = 13298 }
    @{ Row = 41; Old = 1374;  New = 1375 }
    @{ Row = 44; Old = 330;   New = 331 }
    @{ Row = 45; Old = 456;   New = 459 }
    @{ Row = 46; Old = 4075;  New = 4077 }
    @{ Row = 47; Old = 29;    New = 31 }
)

foreach ($u in $updates4) {
    $cell = $ws4.Range("F" + $u.Row)
    $actual = $cell.Value2
    if ($actual -ne $u.Old) {
        Write-Output ("WARNING: 全部类型 F" + $u.Row + " expected " + $u.Old + " but found " + $actual)
    }
    $cell.Value = $u.New
}

$wb.Save()
